$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 37.6319425138707
$ws.Range("C2").Value = 0.7360327430743545

$ws.Range("B3").Value = 45.57433487359113
$ws.Range("C3").Value = 1.340339211219409

$ws.Range("B4").Value = 14.51641272179008
$ws.Range("C4").Value = 1.840780798059799

$ws.Range("B7").Value = 0.0978797837110706
$ws.Range("C7").Value = 0.08012049418058301
